# "Generate Report for Handback" — append a new handback row for the file
# 651a97ba-0481-4a23-abd7-26fd20343089.md to the Overview, zh-cn and de-de
# sheets/tables (new row 4 in each).

$wb = $excel.ActiveWorkbook

$fileName   = "651a97ba-0481-4a23-abd7-26fd20343089.md"
$pathName   = "e2e\651a97ba-0481-4a23-abd7-26fd20343089.md"
$ext        = ".md"
$statusText = "Handed back: in sync with en-US"

$zhXlf      = "651a97ba-0481-4a23-abd7-26fd20343089.d893d3d51ebd27f87f052bc3caca716e8acbf9ec.zh-cn.xlf"
$deXlf      = "651a97ba-0481-4a23-abd7-26fd20343089.d893d3d51ebd27f87f052bc3caca716e8acbf9ec.de-de.xlf"

$handoffZh  = "2016-11-09 01:21:08"
$handoffDe  = "2016-11-09 01:21:22"
$handbackZh = "2016-11-09 01:22:03"
$handbackDe = "2016-11-09 01:22:21"

$srcUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b25190e86b7b7746f059e5d826b1d27954a1ddfa/e2e/651a97ba-0481-4a23-abd7-26fd20343089.md"
$zhTgtUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d893d3d51ebd27f87f052bc3caca716e8acbf9ec/e2e/651a97ba-0481-4a23-abd7-26fd20343089.md"
$deTgtUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d893d3d51ebd27f87f052bc3caca716e8acbf9ec/e2e/651a97ba-0481-4a23-abd7-26fd20343089.md"

# ---------------------------------------------------------------------------
# Overview sheet -> new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusText
$wsOverview.Range("F4").Value = $statusText
$wsOverview.Range("G4").Value = $handoffDe
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcUrl, "", "", $pathName) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# zh-cn sheet -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusText
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $handoffZh
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $handbackZh
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# Leading apostrophe forces text storage so "True"/"False"/"" are written as
# plain shared-string text (matching the source data) instead of being
# auto-coerced to native boolean cells.
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhTgtUrl, "", "", $fileName) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# de-de sheet -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusText
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $handoffDe
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $handbackDe
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# Leading apostrophe forces text storage so "True"/"False"/"" are written as
# plain shared-string text (matching the source data) instead of being
# auto-coerced to native boolean cells.
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deTgtUrl, "", "", $fileName) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))

Write-Output "Handback row appended to Overview, zh-cn and de-de sheets."
